$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Igf2", "Insr", "ECs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 12.840326, 38.520978, 0.3393128690704512, 0.3393128690704511, 109.7033710729426, 987.3303396564839, 0.00865133326887262, 0.008651333268872616),
    @("ECs", "Igf2", "Insr", "FAPs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 9.834223333333334, 29.50267, 0.2598749077175229, 0.2598749077175228, 84.02025396791777, 756.18228571126, 0.006625933289948405, 0.006625933289948403),
    @("ECs", "Igf2", "Insr", "M2", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 8.654269333333334, 25.962808, 0.228693956617749, 0.2286939566177489, 73.93912896291377, 665.4521606662239, 0.005830924246101752, 0.00583092424610175),
    @("ECs", "Igf2", "Insr", "sCs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 6.513324000000001, 19.539972, 0.1721182665942771, 0.1721182665942771, 55.647621383624, 500.828592452616, 0.004388435045352157, 0.004388435045352156),
    @("FAPs", "Igf2", "Insr", "ECs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 12.840326, 38.520978, 0.3393128690704512, 0.3393128690704511, 3496.779228900833, 31471.0130601075, 0.2757600079288161, 0.275760007928816),
    @("FAPs", "Igf2", "Insr", "FAPs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 9.834223333333334, 29.50267, 0.2598749077175229, 0.2598749077175228, 2678.133552401389, 24103.2019716125, 0.211200673906079, 0.211200673906079),
    @("FAPs", "Igf2", "Insr", "M2", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 8.654269333333334, 25.962808, 0.228693956617749, 0.2286939566177489, 2356.799137818889, 21211.19224037, 0.1858598745840339, 0.1858598745840338),
    @("FAPs", "Igf2", "Insr", "sCs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 6.513324000000001, 19.539972, 0.1721182665942771, 0.1721182665942771, 1773.759955495, 15963.839599455, 0.1398807380656027, 0.1398807380656027),
    @("M2", "Igf2", "Insr", "ECs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 12.840326, 38.520978, 0.3393128690704512, 0.3393128690704511, 34.10812176303666, 306.97309586733, 0.002689805478731513, 0.002689805478731512),
    @("M2", "Igf2", "Insr", "FAPs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 9.834223333333334, 29.50267, 0.2598749077175229, 0.2598749077175228, 26.12292607666111, 235.10633468995, 0.002060083817269848, 0.002060083817269848),
    @("M2", "Igf2", "Insr", "M2", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 8.654269333333334, 25.962808, 0.228693956617749, 0.2286939566177489, 22.98858083443111, 206.89722750988, 0.001812905767907927, 0.001812905767907926),
    @("M2", "Igf2", "Insr", "sCs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 6.513324000000001, 19.539972, 0.1721182665942771, 0.1721182665942771, 17.30152708538, 155.71374376842, 0.001364418207135353, 0.001364418207135352),
    @("sCs", "Igf2", "Insr", "ECs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 12.840326, 38.520978, 0.3393128690704512, 0.3393128690704511, 662.0715880589654, 5958.644292530688, 0.05221172239403103, 0.052211722394031),
    @("sCs", "Igf2", "Insr", "FAPs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 9.834223333333334, 29.50267, 0.2598749077175229, 0.2598749077175228, 507.0712269787023, 4563.641042808321, 0.03998821670422561, 0.03998821670422559),
    @("sCs", "Igf2", "Insr", "M2", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 8.654269333333334, 25.962808, 0.228693956617749, 0.2286939566177489, 446.2305583993743, 4016.075025594368, 0.03519025201970542, 0.0351902520197054),
    @("sCs", "Igf2", "Insr", "sCs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 6.513324000000001, 19.539972, 0.1721182665942771, 0.1721182665942771, 335.8393520711681, 3022.554168640513, 0.0264846752761869, 0.02648467527618688),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value2 = $rowVals[$j]
    }
}

Write-Host "Updated rows 2-17"